$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.163.58"
$ws.Range("E2").Value = "  +0.39%  "
$ws.Range("D3").Value = "1.834.73"
$ws.Range("E3").Value = "  +0.32%  "
$ws.Range("D4").Value = "'0.9992"
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'241.86"
$ws.Range("E5").Value = "  +1.25%  "
$ws.Range("D6").Value = "'0.6586"
$ws.Range("E6").Value = "  -0.15%  "
$ws.Range("D7").Value = "'1.000"
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").Value = "'0.07418"
$ws.Range("E8").Value = "  +1.50%  "
$ws.Range("D9").Value = "'0.2931"
$ws.Range("E9").Value = "  -0.30%  "
$ws.Range("D10").Value = "'22.90"
$ws.Range("D11").Value = "'0.07771"
$ws.Range("D12").Value = "1.883.98"
$ws.Range("E12").Value = "  +2.75%  "
$ws.Range("D13").Value = "'4.980"
$ws.Range("E13").Value = "  -0.42%  "
$ws.Range("D14").Value = "'0.6649"
$ws.Range("E14").Value = "  -0.85%  "
$ws.Range("D15").Value = "'82.77"
$ws.Range("E15").Value = "  -3.69%  "
$ws.Range("D16").Value = "'6.106"
$ws.Range("E16").Value = "  +0.10%  "
$ws.Range("D17").Value = "'0.000008558"
$ws.Range("E17").Value = "  +4.60%  "
$ws.Range("D18").Value = "29.177.86"
$ws.Range("E18").Value = "  +0.43%  "
$ws.Range("D19").Value = "2.119.41"
$ws.Range("E19").Value = "  +1.57%  "
$ws.Range("D20").Value = "'226.73"
$ws.Range("E20").Value = "  -0.19%  "
$ws.Range("D21").Value = "'12.44"
$ws.Range("E21").Value = "  +0.11%  "
$ws.Range("D22").Value = "'1.001"
$ws.Range("E22").Value = "  +0.10%  "
$ws.Range("D23").Value = "'7.098"
$ws.Range("E23").Value = "  -1.94%  "
$ws.Range("E24").Value = "  +0.00%  "
$ws.Range("D25").Value = "'159.38"
$ws.Range("E25").Value = "  -0.80%  "
$ws.Range("D26").Value = "'8.596"
$ws.Range("E26").Value = "  -0.54%  "
$ws.Range("D27").Value = "'0.1393"
$ws.Range("E27").Value = "  -1.69%  "
$ws.Range("D28").Value = "'17.92"
$ws.Range("E28").Value = "  +0.12%  "
$ws.Range("D29").Value = "'1.517"
$ws.Range("E29").Value = "  +1.53%  "
$ws.Range("D30").Value = "'4.109"
$ws.Range("E30").Value = "  -2.42%  "
$ws.Range("D31").Value = "'4.041"
$ws.Range("E31").Value = "  -1.29%  "
$ws.Range("D32").Value = "'1.192"
$ws.Range("E32").Value = "  -0.35%  "
$ws.Range("D33").Value = "'0.05270"
$ws.Range("E33").Value = "  -0.32%  "
$ws.Range("D34").Value = "'1.865"
$ws.Range("E34").Value = "  +1.16%  "
$ws.Range("D35").Value = "'0.7384"
$ws.Range("E35").Value = "  -1.44%  "
$ws.Range("D36").Value = "'1.144"
$ws.Range("E36").Value = "  +1.77%  "
$ws.Range("D37").Value = "'2.658"
$ws.Range("E37").Value = "  -0.81%  "
$ws.Range("D38").Value = "1.300.28"
$ws.Range("E38").Value = "  +0.78%  "
$ws.Range("D39").Value = "'0.01793"
$ws.Range("E39").Value = "  -0.43%  "
$ws.Range("D40").Value = "'2.735"
$ws.Range("E40").Value = "  +1.21%  "
$ws.Range("D41").Value = "'0.9191"
$ws.Range("E41").Value = "  -0.17%  "
$ws.Range("D42").Value = "'6.036"
$ws.Range("E42").Value = "  +1.21%  "
$ws.Range("D43").Value = "'0.08589"
$ws.Range("E43").Value = "  +13.84%  "
$ws.Range("E44").Value = "  +0.07%  "
$ws.Range("D45").Value = "'102.66"
$ws.Range("E45").Value = "  -0.70%  "
$ws.Range("D46").Value = "2.034.00"
$ws.Range("E46").Value = "  +2.59%  "
$ws.Range("D47").Value = "'0.5143"
$ws.Range("E47").Value = "  -0.63%  "
$ws.Range("E48").Value = "  -2.19%  "
$ws.Range("D49").Value = "'63.49"
$ws.Range("E49").Value = "  +0.54%  "
$ws.Range("E50").Value = "  +0.33%  "
$ws.Range("D51").Value = "'0.05846"
